$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# D-column values are forced to Text via NumberFormat "@" so Excel
# does not auto-convert decimal-looking strings (e.g. "291.21") into
# numbers; the style is then reset back to Normal so no stray cell
# style index is left behind (matches source cells, which carry no "s").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.456.09'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.31%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.572.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.07%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("E5").Value = '  +0.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3755'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.85'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.66%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3413'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.12%  '

$ws.Range("E10").Value = '  -1.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07571'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.36%  '

$ws.Range("E12").Value = '  -0.02%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.38'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.07%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.022'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.959'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.574.35'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001122'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06739'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.254'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.40'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.90%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.459.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.351'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.69%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.605'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.98%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '148.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.992'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.83%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.749.20'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.028'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.160'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.989'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.875'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08449'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.81%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.379'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02467'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2290'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.55%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06574'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.494'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.38%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.38'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6306'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.819'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.99%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5877'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.103'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.225'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.97%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07333'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.14%  '

